# Updates cryptos list figures (price + 1h volume change) scraped on
# Mon Sep 23 19:37:29 UTC 2024, plus a couple of ranking swaps (rows 29/30
# and 47/48) where two coins traded places.
#
# Price-like text in column D (e.g. "614.57") would be auto-parsed as a
# number by a plain .Value assignment, which also triggers a General->Text
# NumberFormat elsewhere, so for those cells we prefix with a literal
# apostrophe to force text entry, then reset the cell style back to Normal
# so no stray per-cell formatting sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.397.27'
$ws.Range('E2').Value = '  +0.82%  '

# Row 3
$ws.Range('D3').Value = '2.676.25'
$ws.Range('E3').Value = '  +4.23%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').Value = '''614.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.09%  '

# Row 6
$ws.Range('D6').Value = '''143.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.15%  '

# Row 7
$ws.Range('E7').Value = '  +0.14%  '

# Row 8
$ws.Range('E8').Value = '  -0.13%  '

# Row 9
$ws.Range('D9').Value = '2.674.15'
$ws.Range('E9').Value = '  +4.17%  '

# Row 10
$ws.Range('E10').Value = '  +1.10%  '

# Row 11
$ws.Range('D11').Value = '''5.62'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.95%  '

# Row 12
$ws.Range('E12').Value = '  +0.76%  '

# Row 13
$ws.Range('D13').Value = '''0.362'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.98%  '

# Row 14
$ws.Range('E14').Value = '  +1.70%  '

# Row 15
$ws.Range('D15').Value = '3.160.13'
$ws.Range('E15').Value = '  +4.45%  '

# Row 16
$ws.Range('D16').Value = '63.268.79'
$ws.Range('E16').Value = '  +0.79%  '

# Row 17
$ws.Range('E17').Value = '  +0.55%  '

# Row 18
$ws.Range('D18').Value = '2.675.58'
$ws.Range('E18').Value = '  +3.57%  '

# Row 19
$ws.Range('D19').Value = '''11.45'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.95%  '

# Row 20
$ws.Range('D20').Value = '''342.74'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.91%  '

# Row 21
$ws.Range('D21').Value = '''4.41'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.22%  '

# Row 22
$ws.Range('D22').Value = '''6.87'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.94%  '

# Row 24
$ws.Range('D24').Value = '''67.33'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.23%  '

# Row 25
$ws.Range('E25').Value = '  +4.14%  '

# Row 26
$ws.Range('E26').Value = '  -3.04%  '

# Row 27
$ws.Range('D27').Value = '''8.67'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.64%  '

# Row 28
$ws.Range('D28').Value = '''0.164'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.02%  '

# Row 29
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '''1.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.01%  '

# Row 30
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').Value = '''538.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +17.04%  '

# Row 31
$ws.Range('D31').Value = '''7.92'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.63%  '

# Row 32
$ws.Range('E32').Value = '  +7.61%  '

# Row 33
$ws.Range('E33').Value = '  +8.90%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0808'
$ws.Range('E34').Value = '  +1.77%  '

# Row 35
$ws.Range('D35').Value = '''172.34'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.32%  '

# Row 36
$ws.Range('D36').Value = '''5.18'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +14.56%  '

# Row 37
$ws.Range('E37').Value = '  +2.53%  '

# Row 38
$ws.Range('E38').Value = '  +0.00%  '

# Row 39
$ws.Range('D39').Value = '''19.25'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.46%  '

# Row 40
$ws.Range('E40').Value = '  +10.63%  '

# Row 41
$ws.Range('D41').Value = '''176.80'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +12.44%  '

# Row 42
$ws.Range('E42').Value = '  +0.00%  '

# Row 43
$ws.Range('E43').Value = '  +2.23%  '

# Row 44
$ws.Range('D44').Value = '''22.27'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.30%  '

# Row 45
$ws.Range('D45').Value = '''0.0572'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +7.33%  '

# Row 46
$ws.Range('E46').Value = '  +1.15%  '

# Row 47
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.0965'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.80%  '

# Row 48
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '''0.0241'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.62%  '

# Row 49
$ws.Range('D49').Value = '''18.91'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.25%  '

# Row 50
$ws.Range('D50').Value = '''1.75'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +4.99%  '

# Row 51
$ws.Range('E51').Value = '  -0.80%  '
